# Fruta / hortaliza, semanal
# A new weekly price observation is inserted as row 29 ("Vega Modelo de
# Temuco" - Papaya, fecha 2021-09-28, Volumen 50), pushing the previously
# existing rows 29-56 down to rows 30-57.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 29

# Insert a new blank row at position 29; everything at/after it shifts down
# by one (row 56 -> row 57), and the new row inherits the formatting
# (including the date number format on column D) from the row above it.
$ws.Rows.Item($newRow).Insert()

# Seed the new row with the same reference data as the (now shifted down)
# row that used to occupy this slot, then overwrite the two cells that
# actually carry new information for this week's entry.
$sourceRow = $newRow + 1
for ($col = 1; $col -le 20; $col++) {
    $ws.Cells.Item($newRow, $col).Value2 = $ws.Cells.Item($sourceRow, $col).Value2
}

# Fecha (D) and Volumen (M) for the new weekly entry.
$ws.Cells.Item($newRow, 4).Value2 = 44467
$ws.Cells.Item($newRow, 13).Value2 = 50
